# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values on the active worksheet for the listed rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 4
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 0
    18 = 3
    21 = 1
    24 = 1
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
